$d = $word.ActiveDocument

$pairs = @(
    @("668×6=4008", "430×6=2580"),
    @("189×5=945", "891×9=8019"),
    @("259×9=2331", "813×3=2439"),
    @("208×7=1456", "961×7=6727"),
    @("529×9=4761", "878×8=7024"),
    @("914×4=3656", "302×2=604"),
    @("858×4=3432", "109×9=981"),
    @("485×6=2910", "783×3=2349"),
    @("161×4=644", "689×9=6201"),
    @("541×5=2705", "582×4=2328"),
    @("887×9=7983", "905×8=7240"),
    @("729×6=4374", "174×2=348"),
    @("695×6=4170", "293×7=2051"),
    @("943×7=6601", "970×8=7760"),
    @("874×6=5244", "408×3=1224"),
    @("179×8=1432", "504×9=4536"),
    @("274×3=822", "657×4=2628"),
    @("883×6=5298", "992×9=8928"),
    @("829×3=2487", "349×3=1047"),
    @("965×7=6755", "798×7=5586"),
    @("992×5=4960", "430×3=1290"),
    @("182×2=364", "742×7=5194"),
    @("712×5=3560", "422×5=2110"),
    @("954×4=3816", "153×3=459"),
    @("566×6=3396", "157×3=471")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done replacing $($pairs.Count) values"
